$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D1").Value = "required"

$ws.Range("D4").Value = "yes"
$ws.Range("D8").Value = "yes"
$ws.Range("D9").Value = "yes"
$ws.Range("D10").Value = "yes"
$ws.Range("D11").Value = "yes"
$ws.Range("D12").Value = "yes"

$ws.Range("D15").Select()
